$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B (Coin name) swaps ---
$ws.Range('B41').Value = 'Quant'
$ws.Range('B42').Value = 'PaxDollar'
$ws.Range('B47').Value = 'Algorand'
$ws.Range('B48').Value = 'Cronos'

# --- Column C (Link) swaps ---
$ws.Range('C41').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('C42').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('C47').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'

# --- Column D (Price) updates; use leading apostrophe to force text so
#     Excel does not reinterpret these as numbers/dates and lose formatting,
#     then reset the style so no stray quote-prefix style is left behind. ---
$ws.Range('D2').Value = "'25.714.46"
$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').Value = "'1.746.60"
$ws.Range('D3').Style = 'Normal'
$ws.Range('D4').Value = "'1.000"
$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').Value = "'238.48"
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Value = "'0.9992"
$ws.Range('D6').Style = 'Normal'
$ws.Range('D7').Value = "'0.4941"
$ws.Range('D7').Style = 'Normal'
$ws.Range('D8').Value = "'41.53"
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').Value = "'0.2486"
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').Value = "'0.05980"
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').Value = "'1.746.87"
$ws.Range('D11').Style = 'Normal'
$ws.Range('D12').Value = "'0.06778"
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').Value = "'14.85"
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').Value = "'4.476"
$ws.Range('D14').Style = 'Normal'
$ws.Range('D16').Value = "'0.5826"
$ws.Range('D16').Style = 'Normal'
$ws.Range('D17').Value = "'1.000"
$ws.Range('D17').Style = 'Normal'
$ws.Range('D18').Value = "'0.9991"
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').Value = "'25.758.49"
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').Value = "'11.60"
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').Value = "'0.000006567"
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').Value = "'1.965.18"
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').Value = "'3.988"
$ws.Range('D23').Style = 'Normal'
$ws.Range('D24').Value = "'7.925"
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').Value = "'5.041"
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').Value = "'136.49"
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').Value = "'1.487"
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').Value = "'1.836"
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').Value = "'14.57"
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').Value = "'101.00"
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').Value = "'3.801"
$ws.Range('D31').Style = 'Normal'
$ws.Range('D32').Value = "'0.08121"
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').Value = "'3.353"
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').Value = "'0.04420"
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').Value = "'0.9987"
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').Value = "'2.664"
$ws.Range('D36').Style = 'Normal'
$ws.Range('D37').Value = "'1.017"
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').Value = "'0.6070"
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').Value = "'2.699"
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').Value = "'2.061"
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').Value = "'103.82"
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').Value = "'0.9991"
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').Value = "'0.01501"
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').Value = "'0.7761"
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').Value = "'5.186"
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').Value = "'0.3771"
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').Value = "'0.1082"
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').Value = "'0.05126"
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').Value = "'5.981"
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').Value = "'30.42"
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').Value = "'52.67"
$ws.Range('D51').Style = 'Normal'

# --- Column E (Volume/%) updates; same text-safety treatment. ---
$ws.Range('E2').Value = "'  -3.44%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('E3').Value = "'  -5.27%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'  +0.05%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('E5').Value = "'  -8.51%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('E6').Value = "'  -0.06%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = "'  -6.50%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('E8').Value = "'  -7.65%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('E9').Value = "'  -21.78%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('E10').Value = "'  -12.15%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('E11').Value = "'  -5.05%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('E12').Value = "'  -13.04%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('E13').Value = "'  -22.88%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('E15').Value = "'  -12.71%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('E16').Value = "'  -25.91%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('E17').Value = "'  +0.06%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('E18').Value = "'  -0.12%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('E19').Value = "'  -3.37%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('E20').Value = "'  -16.70%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('E21').Value = "'  -17.27%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('E22').Value = "'  -5.21%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('E23').Value = "'  -13.68%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('E24').Value = "'  -15.22%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('E25').Value = "'  -16.06%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('E26').Value = "'  -4.54%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('E27').Value = "'  -12.79%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('E28').Value = "'  -17.33%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('E29').Value = "'  -14.62%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('E30').Value = "'  -9.25%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('E31').Value = "'  -10.03%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('E32').Value = "'  -6.79%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('E33').Value = "'  -17.99%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('E34').Value = "'  -9.12%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('E35').Value = "'  -0.05%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('E36').Value = "'  -7.12%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('E37').Value = "'  -10.84%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('E38').Value = "'  -17.12%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('E39').Value = "'  -12.96%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('E40').Value = "'  -12.42%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('E41').Value = "'  -5.11%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('E42').Value = "'  -0.06%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('E43').Value = "'  -13.65%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('E44').Value = "'  -14.14%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('E45').Value = "'  -12.33%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('E46').Value = "'  -22.16%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('E47').Value = "'  -13.34%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('E48').Value = "'  -11.95%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('E49').Value = "'  -22.41%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('E50').Value = "'  -12.97%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('E51').Value = "'  -12.42%  "
$ws.Range('E51').Style = 'Normal'
